# Update Leave Card 8/7/2023 4:34 PM
$wb = $excel.ActiveWorkbook

# --- workbook.xml: absPath url (cosmetic, last-saved machine path) ---
# Not exposed via a documented COM property; skip (no functional effect).

$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row within the leave-record table at worksheet row 561 ---
# (shifts existing rows 561..780 down to 562..781)
$ws.Rows.Item(561).Insert()

# Copy formatting from the row directly below (which now holds the record
# that used to be at row 561) so the new blank row matches the table's
# normal row styling.
$ws.Range("A562:K562").Copy()
$ws.Range("A561:K561").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the newly inserted row with the new leave entry ---
$ws.Cells.Item(561, 1).Value = $null                     # A561 PERIOD - blank
$ws.Cells.Item(561, 2).Value = "A(3-0-0)"                 # B561 PARTICULARS
$ws.Cells.Item(561, 3).Value = $null                      # C561 EARNED - blank
$ws.Cells.Item(561, 4).Value = 3                          # D561 Absence Undertime W/ Pay
$ws.Cells.Item(561, 5).Value = $null                      # E561 BALANCE - blank
$ws.Cells.Item(561, 6).Value = $null                      # F561 Absence Undertime W/O Pay - blank
$ws.Cells.Item(561, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Cells.Item(561, 8).Value = $null                      # H561 Absence Undertime  W/ Pay - blank
$ws.Cells.Item(561, 9).Value = $null                      # I561 BALANCE  - blank
$ws.Cells.Item(561, 10).Value = $null                     # J561 Absence Undertime  W/O Pay - blank

# K561 uses the plain left-aligned text style (matching other REMARKS-only
# rows), rather than the highlighted style copied from the row below.
$ws.Range("K567").Copy()
$ws.Range("K561").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(561, 11).Value = "5/5,6,10/2022"           # K561 REMARKS

# --- Resize Table1 so it spans the newly inserted row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K781"))

# Re-assert the calculated-column formula on the table's last (totals-style)
# row; the recalc triggered mid-insert can otherwise rewrite it into a
# transient structured-reference form that evaluates to #VALUE!.
$ws.Cells.Item(781, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Restore the active selection to match the authored workbook ---
$ws.Range("K561").Select() | Out-Null

$wb.Application.Calculate()
